# Insert a new weekly price group (date 2021-09-16 / serial 44455) for
# "Comercializadora del Agro de Limarí - Pepino dulce" ahead of the
# existing 2021-08-19 group at row 185. This pushes every subsequent row
# down by 4 (the sheet grows from A1:R252 to A1:R256) and the last
# existing group gets a matching new weekly row appended at the tail
# automatically as part of the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 185..188, shifting everything below down.
$ws.Rows("185:188").Insert()

# Common (unchanging) field values shared by every data row in this sheet.
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fecha = 44455

# New group rows: Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$newRows = @(
    @{ Row = 185; Calidad = "Especial"; Volumen = 240; PMin = 13000; PMax = 13500; PProm = 13250; PKg = 736 },
    @{ Row = 186; Calidad = "Primera";  Volumen = 360; PMin = 11000; PMax = 11500; PProm = 11250; PKg = 625 },
    @{ Row = 187; Calidad = "Segunda";  Volumen = 360; PMin = 9000;  PMax = 9500;  PProm = 9250;  PKg = 514 },
    @{ Row = 188; Calidad = "Tercera";  Volumen = 240; PMin = 6000;  PMax = 6500;  PProm = 6250;  PKg = 347 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
